# Apply updates described in commit message:
# "updated endpoint /api/user/login for all users
#  changed the endpoint to /api/user/login and tested. Test suit: test/user-login.test.js"

$wb = $excel.ActiveWorkbook

# --- Sheet "routes": update status / remarks for the login endpoint, and
#     reword the "register" sub-route description.
$wsRoutes = $wb.Worksheets.Item("routes")
$wsRoutes.Activate() | Out-Null
$wsRoutes.Range("F2").Value = "done"
$wsRoutes.Range("H2").Value = "end point is /api/user/login now"
$wsRoutes.Range("D3").Value = "registration for course"
$wsRoutes.Range("H2").Select() | Out-Null

# --- Sheet "tasks": record the work done / tested, and add a new task row.
$wsTasks = $wb.Worksheets.Item("tasks")
$wsTasks.Activate() | Out-Null
$wsTasks.Range("D2").Value = "done"
$wsTasks.Range("E2").Value = "tested with jest and supertest. Test suit: test/user-login.test.js"

$wsTasks.Range("A3").Value = 45554
$wsTasks.Range("A3").NumberFormat = $wsTasks.Range("A2").NumberFormat
$wsTasks.Range("B3").Value = "Jasdeep"
$wsTasks.Range("C3").Value = "finalize student schema - merging student and applicant"

$wsTasks.Range("D3").Select() | Out-Null

$wsRoutes.Activate() | Out-Null
